$wb = $excel.ActiveWorkbook

# --- 1) Insert the new "1M_Statistics" sheet ---
#     Rename the existing "Annual_Statistics" sheet to "1M_Statistics" in place
#     (this keeps its original identity/sheetId and all of its header/
#     label/style/column-width formatting intact), then add a fresh copy right
#     after it and rename that copy back to "Annual_Statistics". That mirrors
#     how the workbook really gained the new sheet: a brand new sheet slot was
#     created for "Annual_Statistics" while "1M_Statistics" reuses the old slot.
$oneMonth = $wb.Worksheets.Item("Annual_Statistics")
$oneMonth.Name = "1M_Statistics"
$oneMonth.Copy($null, $oneMonth)
$annual = $wb.Worksheets.Item("1M_Statistics (2)")
$annual.Name = "Annual_Statistics"

# --- 2) Update the new "1M_Statistics" sheet values ---
$oneMonth.Range("F2").Value = 191.3809661865234
$oneMonth.Range("B3").Value = 0.03465597060360763
$oneMonth.Range("C3").Value = 0.008833839573015036
$oneMonth.Range("D3").Value = 0.009760939124285642
$oneMonth.Range("E3").Value = 0.03776925070725742
$oneMonth.Range("F3").Value = 0.03685140124486985
$oneMonth.Range("B4").Value = 0.05547189063155183
$oneMonth.Range("C4").Value = 0.04012324282203933
$oneMonth.Range("D4").Value = 0.04457046772956886
$oneMonth.Range("E4").Value = 0.1198100025669567
$oneMonth.Range("F4").Value = 0.05751045451745348
$oneMonth.Range("B5").Value = 0.003077130650238848
$oneMonth.Range("C5").Value = 0.001609874614556331
$oneMonth.Range("D5").Value = 0.001986526593632539
$oneMonth.Range("E5").Value = 0.01435443671509417
$oneMonth.Range("F5").Value = 0.003307452378804086
$oneMonth.Range("F6").Value = -0.05248821566458358
$oneMonth.Range("F7").Value = 1.470635920535974
$oneMonth.Range("B8").Value = -0.08294049223427413
$oneMonth.Range("C8").Value = -0.0607929951813917
$oneMonth.Range("D8").Value = -0.06531872274428896
$oneMonth.Range("E8").Value = -0.1663489320221079
$oneMonth.Range("F8").Value = -0.07632621153963082
$oneMonth.Range("B9").Value = -0.09586675600067596
$oneMonth.Range("C9").Value = -0.08173464394712393
$oneMonth.Range("D9").Value = -0.0929467329417244
$oneMonth.Range("E9").Value = -0.2019226716669809
$oneMonth.Range("F9").Value = -0.1096430311662466
$oneMonth.Range("B10").Value = -0.1063411838753529
$oneMonth.Range("C10").Value = -0.0907415442651075
$oneMonth.Range("D10").Value = -0.1057585957915729
$oneMonth.Range("E10").Value = -0.2165793652403061
$oneMonth.Range("F10").Value = -0.1220993020676024
$oneMonth.Range("B11").Value = -0.1145941197849433
$oneMonth.Range("C11").Value = -0.1208543551483093
$oneMonth.Range("D11").Value = -0.1273253158243038
$oneMonth.Range("E11").Value = -0.2336711219351461
$oneMonth.Range("F11").Value = -0.1615451349056326
$oneMonth.Range("F12").Value = 0.5431446161467421
$oneMonth.Range("F14").Value = 3.637364727037048

# --- 3) Update the "Annual_Statistics" sheet values (tiny recompute deltas) ---
$annual.Range("F2").Value = 191.3809661865234
$annual.Range("F3").Value = 0.438707157677022
$annual.Range("F4").Value = 0.1984299198230978
$annual.Range("F5").Value = 0.03937443308100102
$annual.Range("F6").Value = -0.05248821566458358
$annual.Range("F7").Value = 1.470635920535974
$annual.Range("B8").Value = -0.2582473431006131
$annual.Range("C8").Value = -0.1945888538675579
$annual.Range("D8").Value = -0.2079007388543198
$annual.Range("E8").Value = -0.4662105678841544
$annual.Range("F8").Value = -0.2396248156230657
$annual.Range("B9").Value = -0.2937026700115168
$annual.Range("C9").Value = -0.2548766893648746
$annual.Range("D9").Value = -0.2858009892613433
$annual.Range("E9").Value = -0.5407776151820956
$annual.Range("F9").Value = -0.3301465057954592
$annual.Range("B10").Value = -0.3215364639872514
$annual.Range("C10").Value = -0.2797922171780177
$annual.Range("D10").Value = -0.3200091667584076
$annual.Range("E10").Value = -0.5692274471813069
$annual.Range("F10").Value = -0.3619303733051077
$annual.Range("B11").Value = -0.3429113583944799
$annual.Range("C11").Value = -0.358802938435328
$annual.Range("D11").Value = -0.3749405665416327
$annual.Range("E11").Value = -0.6007963321270422
$annual.Range("F11").Value = -0.4555226173149048
$annual.Range("F12").Value = 0.5431446161467421
$annual.Range("F14").Value = 3.637364727037048

# --- 4) Update "1D_Statistics" values ---
$s1 = $wb.Worksheets.Item("1D_Statistics")
$s1.Range("F2").Value = 191.3809661865234
$s1.Range("F3").Value = 0.001754828630708088
$s1.Range("F4").Value = 0.0125498100513117
$s1.Range("F5").Value = 0.0001574977323240041
$s1.Range("F6").Value = -0.05248821566458358
$s1.Range("F7").Value = 1.470635920535974
$s1.Range("F9").Value = -0.02534225636928439
$s1.Range("F10").Value = -0.02841672472349273
$s1.Range("F11").Value = -0.03844879802287123
$s1.Range("F12").Value = 0.5431446161467421
$s1.Range("F14").Value = 3.637364727037048

# --- 5) Update "3D_Statistics" values ---
$s3 = $wb.Worksheets.Item("3D_Statistics")
$s3.Range("F2").Value = 191.3809661865234
$s3.Range("F3").Value = 0.005264485892124264
$s3.Range("F4").Value = 0.02173690863421044
$s3.Range("F5").Value = 0.0004724931969720122
$s3.Range("F6").Value = -0.05248821566458358
$s3.Range("F7").Value = 1.470635920535974
$s3.Range("F9").Value = -0.04389407561003655
$s3.Range("F10").Value = -0.04921921100578806
$s3.Range("F11").Value = -0.06659527166556675
$s3.Range("F12").Value = 0.5431446161467421
$s3.Range("F14").Value = 3.637364727037048

# --- 6) Update "5D_Statistics" values ---
$s5 = $wb.Worksheets.Item("5D_Statistics")
$s5.Range("F2").Value = 191.3809661865234
$s5.Range("F3").Value = 0.008774143153540439
$s5.Range("F4").Value = 0.02806222837944308
$s5.Range("F5").Value = 0.0007874886616200204
$s5.Range("F6").Value = -0.05248821566458358
$s5.Range("F7").Value = 1.470635920535974
$s5.Range("F9").Value = -0.05666700794494692
$s5.Range("F10").Value = -0.06354172817962866
$s5.Range("F11").Value = -0.08597412603229958
$s5.Range("F12").Value = 0.5431446161467421
$s5.Range("F14").Value = 3.637364727037048

# --- 7) Update "10D_Statistics" values ---
$s10 = $wb.Worksheets.Item("10D_Statistics")
$s10.Range("F2").Value = 191.3809661865234
$s10.Range("F3").Value = 0.01754828630708088
$s10.Range("F4").Value = 0.03968598396461956
$s10.Range("F5").Value = 0.001574977323240041
$s10.Range("F6").Value = -0.05248821566458358
$s10.Range("F7").Value = 1.470635920535974
$s10.Range("B8").Value = -0.05799800594496463
$s10.Range("C8").Value = -0.04235724559895471
$s10.Range("D8").Value = -0.04554397971420299
$s10.Range("E8").Value = -0.1179890145847151
$s10.Range("F8").Value = -0.05331480146745704
$s10.Range("B9").Value = -0.06718072179566781
$s10.Range("C9").Value = -0.05714343352132234
$s10.Range("D9").Value = -0.0651028219970956
$s10.Range("E9").Value = -0.1441362940579037
$s10.Range("F9").Value = -0.07701218969951162
$s10.Range("B10").Value = -0.07465154901957483
$s10.Range("C10").Value = -0.06353497497395366
$s10.Range("D10").Value = -0.07423531109854076
$s10.Range("E10").Value = -0.155013822350441
$s10.Range("F10").Value = -0.0859422939220238
$s10.Range("B11").Value = -0.08055703963017125
$s10.Range("C11").Value = -0.08504801289619279
$s10.Range("D11").Value = -0.08970058675715831
$s10.Range("E11").Value = -0.1677785513268204
$s10.Range("F11").Value = -0.1144849050085754
$s10.Range("F12").Value = 0.5431446161467421
$s10.Range("F14").Value = 3.637364727037048

# --- 8) Restore the originally active sheet/selection ---
$s1.Activate()

